$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.2460971570146455
$ws.Range("E2").Value = 0.1773659978597841
$ws.Range("F2").Value = 1.012249720703679
$ws.Range("G2").Value = 0.5272962014680189
$ws.Range("H2").Value = 0.5812333897015094
$ws.Range("J2").Value = 0.1772330055721199
$ws.Range("O2").Value = 2.183072434228421
$ws.Range("D3").Value = 0.2444676174135338
$ws.Range("E3").Value = 0.1742912681867494
$ws.Range("F3").Value = 0.9831260700589439
$ws.Range("G3").Value = 0.4999711970627061
$ws.Range("H3").Value = 0.5723711798591182
$ws.Range("J3").Value = 0.1721769985553436
$ws.Range("O3").Value = 2.106238768013611
$ws.Range("D4").Value = 0.2435737801816771
$ws.Range("E4").Value = 0.1724976178116577
$ws.Range("F4").Value = 0.965822117505752
$ws.Range("G4").Value = 0.4835087855672953
$ws.Range("H4").Value = 0.5672287518405454
$ws.Range("J4").Value = 0.1691764068401369
$ws.Range("O4").Value = 2.060343181411412
$ws.Range("D5").Value = 0.2432364075528497
$ws.Range("E5").Value = 0.1717903946035548
$ws.Range("F5").Value = 0.9589159623540695
$ws.Range("G5").Value = 0.4768793569706844
$ws.Range("H5").Value = 0.5652084079914204
$ws.Range("J5").Value = 0.1679797192786268
$ws.Range("O5").Value = 2.041962182227792
$ws.Range("D6").Value = 0.2431820113789342
$ws.Range("E6").Value = 0.1716743925400941
$ws.Range("F6").Value = 0.957777978631924
$ws.Range("G6").Value = 0.4757833212342462
$ws.Range("H6").Value = 0.5648774770936456
$ws.Range("J6").Value = 0.1677825844649661
$ws.Range("O6").Value = 2.038929463360006
$ws.Range("D7").Value = 0.2435691213926958
$ws.Range("E7").Value = 0.1724879839688604
$ws.Range("F7").Value = 0.9657283902091223
$ws.Range("G7").Value = 0.483419058529023
$ws.Range("H7").Value = 0.5672012000805466
$ws.Range("J7").Value = 0.1691601623405177
$ws.Range("O7").Value = 2.0600939860322
$ws.Range("D8").Value = 0.2455131683177569
$ws.Range("E8").Value = 0.1762862718188565
$ws.Range("F8").Value = 1.002087880873759
$ws.Range("G8").Value = 0.5178090186325619
$ws.Range("H8").Value = 0.5781156669492873
$ws.Range("J8").Value = 0.1754681284193538
$ws.Range("O8").Value = 2.156314253543428
$ws.Range("D9").Value = 0.2501707657315677
$ws.Range("E9").Value = 0.1844829058625095
$ws.Range("F9").Value = 1.077981530534075
$ws.Range("G9").Value = 0.5877591867548517
$ws.Range("H9").Value = 0.6018910408742215
$ws.Range("J9").Value = 0.1886638921439214
$ws.Range("O9").Value = 2.355184202376734
$ws.Range("D10").Value = 0.2541069149564663
$ws.Range("E10").Value = 0.1909625588303498
$ws.Range("F10").Value = 1.136556074637113
$ws.Range("G10").Value = 0.6407026906573208
$ws.Range("H10").Value = 0.6208072618303788
$ws.Range("J10").Value = 0.1988664766723076
$ws.Range("O10").Value = 2.507549276033956
$ws.Range("D11").Value = 0.2560090365130065
$ws.Range("E11").Value = 0.194010052836056
$ws.Range("F11").Value = 1.163818320499843
$ws.Range("G11").Value = 0.665129526491171
$ws.Range("H11").Value = 0.629727949599868
$ws.Range("J11").Value = 0.2036191076019378
$ws.Range("O11").Value = 2.578233826267081
$ws.Range("D12").Value = 0.2567453311371821
$ws.Range("E12").Value = 0.1951784315297758
$ws.Range("F12").Value = 1.174230611703237
$ws.Range("G12").Value = 0.674428840240239
$ws.Range("H12").Value = 0.633151365407457
$ws.Range("J12").Value = 0.2054348799556607
$ws.Range("O12").Value = 2.605198180176444
$ws.Range("D13").Value = 0.2565860458039424
$ws.Range("E13").Value = 0.1949261616899847
$ws.Range("F13").Value = 1.171984194306717
$ws.Range("G13").Value = 0.6724238672802585
$ws.Range("H13").Value = 0.6324120554357648
$ws.Range("J13").Value = 0.2050431063755838
$ws.Range("O13").Value = 2.599382130016465
$ws.Range("D14").Value = 0.2560692914179441
$ws.Range("E14").Value = 0.1941058882924764
$ws.Range("F14").Value = 1.164673168681375
$ws.Range("G14").Value = 0.6658935957561596
$ws.Range("H14").Value = 0.6300086875166357
$ws.Range("J14").Value = 0.2037681703561276
$ws.Range("O14").Value = 2.58044823741966
$ws.Range("D15").Value = 0.2557548473064344
$ws.Range("E15").Value = 0.1936053170410972
$ws.Range("F15").Value = 1.160206502899541
$ws.Range("G15").Value = 0.6619000535302746
$ws.Range("H15").Value = 0.6285424593910989
$ws.Range("J15").Value = 0.2029893268937144
$ws.Range("O15").Value = 2.568876427918497
$ws.Range("D16").Value = 0.2539848482162341
$ws.Range("E16").Value = 0.19076540734887
$ws.Range("F16").Value = 1.134786829026837
$ws.Range("G16").Value = 0.6391132448457597
$ws.Range("H16").Value = 0.6202306213540965
$ws.Range("J16").Value = 0.1985581258946922
$ws.Range("O16").Value = 2.502957520220775
$ws.Range("D17").Value = 0.2529275540797755
$ws.Range("E17").Value = 0.1890487917632555
$ws.Range("F17").Value = 1.119350579213631
$ws.Range("G17").Value = 0.62522211588049
$ws.Range("H17").Value = 0.6152123840539048
$ws.Range("J17").Value = 0.1958682884033607
$ws.Range("O17").Value = 2.462870201680346
$ws.Range("D18").Value = 0.2523299280996838
$ws.Range("E18").Value = 0.1880708399873043
$ws.Range("F18").Value = 1.110530072250768
$ws.Range("G18").Value = 0.6172645418874367
$ws.Range("H18").Value = 0.6123557350697695
$ws.Range("J18").Value = 0.1943316508988744
$ws.Range("O18").Value = 2.43994224790481
$ws.Range("D19").Value = 0.252129387036959
$ws.Range("E19").Value = 0.1877413368565612
$ws.Range("F19").Value = 1.107553563374211
$ws.Range("G19").Value = 0.6145757763602546
$ws.Range("H19").Value = 0.6113936264415543
$ws.Range("J19").Value = 0.1938131718419669
$ws.Range("O19").Value = 2.432201422091566
$ws.Range("D20").Value = 0.2530390182404858
$ws.Range("E20").Value = 0.1892305555620197
$ws.Range("F20").Value = 1.120987789708067
$ws.Range("G20").Value = 0.6266975129547916
$ws.Range("H20").Value = 0.6157435095531412
$ws.Range("J20").Value = 0.1961535407199193
$ws.Range("O20").Value = 2.467124190622883
$ws.Range("D21").Value = 0.2562206407072125
$ws.Range("E21").Value = 0.1943464327932958
$ws.Range("F21").Value = 1.16681818736771
$ws.Range("G21").Value = 0.6678103529509372
$ws.Range("H21").Value = 0.6307133845864143
$ws.Range("J21").Value = 0.2041422140131885
$ws.Range("O21").Value = 2.586004214704246
$ws.Range("D22").Value = 0.2583932720500854
$ws.Range("E22").Value = 0.1977736478161489
$ws.Range("F22").Value = 1.197287856512375
$ws.Range("G22").Value = 0.6949679831682829
$ws.Range("H22").Value = 0.6407613510683916
$ws.Range("J22").Value = 0.2094568644433963
$ws.Range("O22").Value = 2.664851336729839
$ws.Range("D23").Value = 0.2572251772793663
$ws.Range("E23").Value = 0.1959368213279973
$ws.Range("F23").Value = 1.180978317030494
$ws.Range("G23").Value = 0.6804470454144109
$ws.Range("H23").Value = 0.6353743924378819
$ws.Range("J23").Value = 0.2066117615123773
$ws.Range("O23").Value = 2.622663646788794
$ws.Range("D24").Value = 0.2529885934366689
$ws.Range("E24").Value = 0.1891483522606663
$ws.Range("F24").Value = 1.120247438639609
$ws.Range("G24").Value = 0.6260303967945902
$ws.Range("H24").Value = 0.615503299234291
$ws.Range("J24").Value = 0.1960245476697366
$ws.Range("O24").Value = 2.465200592235419
$ws.Range("D25").Value = 0.248820380564311
$ws.Range("E25").Value = 0.1821852480318924
$ws.Range("F25").Value = 1.056956987943408
$ws.Range("G25").Value = 0.5685648655341993
$ws.Range("H25").Value = 0.595204980215243
$ws.Range("J25").Value = 0.1850052413104351
$ws.Range("O25").Value = 2.300289579330752
